# Fix(module3): use uncon_planned_qty for future production; keep produced for today
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductionPlan")

# Row 2 (MAT_A): uncon_planned_qty / con_planned_qty 1050 -> 860, produced_qty 999 -> 817
$ws.Range("G2").Value = 860
$ws.Range("H2").Value = 860
$ws.Range("J2").Value = 817

# Row 3 (MAT_B): uncon_planned_qty / con_planned_qty 80 -> 112, produced_qty 70 -> 99
$ws.Range("G3").Value = 112
$ws.Range("H3").Value = 112
$ws.Range("J3").Value = 99
